$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
@("ECs","Gm13306","Ccr10","ECs",1,0.5,0.169054,0.338108,0.2998784017861604,0.2245827615387374,1,0.5,0.2651005,0.530201,0.2682690254597046,0.2682690254597046,0.04481629992700001,0.179265199708,0.08044808660358699,0.06024859857304632),
@("ECs","Gm13306","Ccr10","MuSCs",1,0.5,0.169054,0.338108,0.2998784017861604,0.2245827615387374,1,0.5,0.7230885,1.446177,0.7317309745402955,0.7317309745402955,0.122241003279,0.488964013116,0.2194303151825734,0.1643341629656911),
@("FAPs","Gm13306","Ccr10","ECs",3,1,0.3036856666666667,0.911057,0.5386963477076239,0.605154852825717,1,0.5,0.2651005,0.530201,0.2682690254597046,0.2682690254597046,0.08050722207616667,0.483043332457,0.1445155442182264,0.1623443026197661),
@("FAPs","Gm13306","Ccr10","MuSCs",3,1,0.3036856666666667,0.911057,0.5386963477076239,0.605154852825717,1,0.5,0.7230885,1.446177,0.7317309745402955,0.7317309745402955,0.2195916131815,1.317549679089,0.3941808034893975,0.442810550205951),
@("MuSCs","Gm13306","Ccr10","ECs",1,0.5,0.0166775,0.033355,0.02958357711612082,0.02215551838798428,1,0.5,0.2651005,0.530201,0.2682690254597046,0.2682690254597046,0.00442121358875,0.017684854355,0.007936357402553752,0.005943639326499107),
@("MuSCs","Gm13306","Ccr10","MuSCs",1,0.5,0.0166775,0.033355,0.02958357711612082,0.02215551838798428,1,0.5,0.7230885,1.446177,0.7317309745402955,0.7317309745402955,0.01205930845875,0.048237233835,0.02164721971356708,0.01621187906148517),
@("Neutrophils","Gm13306","Ccr10","ECs",2,0.6666666666666666,0.07432466666666666,0.222974,0.1318416733900949,0.1481068672475613,1,0.5,0.2651005,0.530201,0.2682690254597046,0.2682690254597046,0.01970350629566667,0.118221037774,0.03536903723533744,0.0397324849403931),
@("Neutrophils","Gm13306","Ccr10","MuSCs",2,0.6666666666666666,0.07432466666666666,0.222974,0.1318416733900949,0.1481068672475613,1,0.5,0.7230885,1.446177,0.7317309745402955,0.7317309745402955,0.053743311733,0.322459870398,0.0964726361547575,0.1083743823071682)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 2
    $row = $data[$i]
    for ($j = 0; $j -lt $row.Length; $j++) {
        $c = $j + 1
        $ws.Cells.Item($r, $c).Value = $row[$j]
    }
}
